# remove general cat and brand in assortments
#
# On the "Assortment" sheet, row 2 defines an exclusion rule (A2 = scene
# type). This change fills in the "categories_to_exclude" (B2) and
# "brands_to_exclude" (C2) columns with "General" so that the "General"
# category/brand gets excluded ("removed") from that assortment rule.
# The Assortment sheet also becomes the active/selected workbook tab.

$wb = $excel.ActiveWorkbook

$assortment = $wb.Worksheets.Item("Assortment")

# Populate the "categories_to_exclude" / "brands_to_exclude" cells for the
# second data row with "General" (creates/reuses the shared string).
$assortment.Range("B2").Value = "General"
$assortment.Range("C2").Value = "General"

# Update the cell selection on the Assortment sheet to C2 (matches the
# new <selection activeCell="C2" .../> recorded in the workbook).
$assortment.Range("C2").Select() | Out-Null

# Make "Assortment" the active/selected sheet (it becomes the workbook's
# active tab, replacing "BayCountKPI").
$assortment.Activate()
